$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old "total" (row 2) and "avg" (row 3) rows from the top ---
# (they get re-added at the bottom later, after the fund rows)
$ws.Rows("2:3").Delete()

# After the delete, the sheet looks like:
#   Row1: B1 = 06-01-2023
#   Row2: Alpha Acciones            B2=250446.63
#   Row3: Alpha Mega                B3=251985.17
#   Row4: HF Acciones Argentinas    B4=21066.07
#   Row5: HF Acciones Lideres       B5=37199.5

# --- Add the new date column header (cloning the style of the existing date cell) ---
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Prepare three new rows (6,7,8) with the same label styling as the ---
# --- existing fund rows, by cloning the format of an existing styled cell ---
$ws.Range("A5").Copy()
$ws.Range("A6:A8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 6: Toronto trust Argy (note: no historical value for the first date,
# leaving B6 as an explicit empty string since this fund wasn't tracked yet)
$ws.Range("A6").Value = "Toronto trust Argy"
$ws.Range("B6").Formula = '=""'
$ws.Range("C6").Value = 232168.01

# Row 7: avg
$ws.Range("A7").Value = "avg"
$ws.Range("B7").Value = 140174.34
$ws.Range("C7").Value = 156360.61

# Row 8: total
$ws.Range("A8").Value = "total"
$ws.Range("B8").Value = 560697.37
$ws.Range("C8").Value = 781803.0699999999

# --- Fill in column C for the existing fund rows ---
$ws.Range("C2").Value = 244796.51
$ws.Range("C3").Value = 250883.82
$ws.Range("C4").Value = 19393.31
$ws.Range("C5").Value = 34561.42

Write-Output "edit applied"
